# "Tested req 45 to 88" - fill in column I (the tester's results column for
# this requirement-testing group) for requirement rows 5-94 that were tested.
# Values are the free-text result the tester typed: "Unable to Test",
# "Test Working" / "Test working", or "Uable to Test" (typo, as originally
# entered by the tester).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I5').Value = 'Unable to Test'
$ws.Range('I6').Value = 'Unable to Test'
$ws.Range('I11').Value = 'Unable to Test'
$ws.Range('I12').Value = 'Unable to Test'
$ws.Range('I15').Value = 'Unable to Test'
$ws.Range('I16').Value = 'Unable to Test'
$ws.Range('I17').Value = 'Unable to Test'
$ws.Range('I18').Value = 'Unable to Test'
$ws.Range('I19').Value = 'Unable to Test'
$ws.Range('I20').Value = 'Unable to Test'
$ws.Range('I21').Value = 'Unable to Test'
$ws.Range('I22').Value = 'Unable to Test'
$ws.Range('I23').Value = 'Unable to Test'
$ws.Range('I26').Value = 'Unable to Test'
$ws.Range('I27').Value = 'Unable to Test'
$ws.Range('I28').Value = 'Unable to Test'
$ws.Range('I29').Value = 'Unable to Test'
$ws.Range('I31').Value = 'Unable to Test'
$ws.Range('I32').Value = 'Unable to Test'
$ws.Range('I33').Value = 'Unable to Test'
$ws.Range('I35').Value = 'Unable to Test'
$ws.Range('I36').Value = 'Unable to Test'
$ws.Range('I37').Value = 'Unable to Test'
$ws.Range('I38').Value = 'Unable to Test'
$ws.Range('I39').Value = 'Unable to Test'
$ws.Range('I41').Value = 'Unable to Test'
$ws.Range('I42').Value = 'Unable to Test'
$ws.Range('I43').Value = 'Unable to Test'
$ws.Range('I46').Value = 'Unable to Test'
$ws.Range('I47').Value = 'Unable to Test'
$ws.Range('I48').Value = 'Unable to Test'
$ws.Range('I49').Value = 'Unable to Test'
$ws.Range('I50').Value = 'Unable to Test'
$ws.Range('I51').Value = 'Unable to Test'
$ws.Range('I52').Value = 'Unable to Test'
$ws.Range('I53').Value = 'Test Working'
$ws.Range('I54').Value = 'Test working'
$ws.Range('I55').Value = 'Test working'
$ws.Range('I56').Value = 'Test Working'
$ws.Range('I57').Value = 'Test Working'
$ws.Range('I58').Value = 'Unable to Test'
$ws.Range('I59').Value = 'Test Working'
$ws.Range('I60').Value = 'Uable to Test'
$ws.Range('I61').Value = 'Test Working'
$ws.Range('I62').Value = 'Unable to Test'
$ws.Range('I63').Value = 'Unable to Test'
$ws.Range('I64').Value = 'Unable to Test'
$ws.Range('I65').Value = 'Unable to Test'
$ws.Range('I66').Value = 'Unable to Test'
$ws.Range('I67').Value = 'Test Working'
$ws.Range('I68').Value = 'Test Working'
$ws.Range('I69').Value = 'Unable to Test'
$ws.Range('I70').Value = 'Unable to Test'
$ws.Range('I71').Value = 'Unable to Test'
$ws.Range('I72').Value = 'Unable to Test'
$ws.Range('I73').Value = 'Unable to Test'
$ws.Range('I74').Value = 'Unable to Test'
$ws.Range('I75').Value = 'Uable to Test'
$ws.Range('I76').Value = 'Test Working'
$ws.Range('I77').Value = 'Unable to Test'
$ws.Range('I78').Value = 'Unable to Test'
$ws.Range('I79').Value = 'Unable to Test'
$ws.Range('I80').Value = 'Unable to Test'
$ws.Range('I81').Value = 'Unable to Test'
$ws.Range('I82').Value = 'Unable to Test'
$ws.Range('I83').Value = 'Unable to Test'
$ws.Range('I84').Value = 'Test Working'
$ws.Range('I85').Value = 'Test Working'
$ws.Range('I86').Value = 'Uable to Test'
$ws.Range('I87').Value = 'Unable to Test'
$ws.Range('I88').Value = 'Unable to Test'
$ws.Range('I89').Value = 'Test Working'
$ws.Range('I90').Value = 'Uable to Test'
$ws.Range('I91').Value = 'Test Working'
$ws.Range('I92').Value = 'Test Working'
$ws.Range('I93').Value = 'Unable to Test'
$ws.Range('I94').Value = 'Unable to Test'

# Match the author's final on-screen selection after this editing session.
$ws.Range("C89").Select()
